# Added multithread and export of file
# Fill in the newly-collected rows (9 through 15 in the "No." column, i.e.
# worksheet rows 8-15) with the received submissions, and refresh the
# "MAIN MSG" (column J) + country (column K) text for the three rows that
# already had a parsed MAIN MSG (rows 13-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("B8").Value = 45154
$ws.Range("C8").Value = "Madhukar Verma"
$ws.Range("D8").Value = "JKR"
$ws.Range("E8").Value = "ENGLISH"
$ws.Range("F8").Value = "42729 Mayfair Park Ave Fremont Fremont 94538 California USA"
$ws.Range("G8").Value = 6509194563

# Row 9
$ws.Range("B9").Value = 45154
$ws.Range("C9").Value = "Hina Kanjiani"
$ws.Range("D9").Value = "GG"
$ws.Range("E9").Value = "ENGLISH"
$ws.Range("F9").Value = "2050 Cross Creek Ct Allen 75013 TX USA"
$ws.Range("G9").Value = 6825577801

# Row 10
$ws.Range("B10").Value = 45154
$ws.Range("C10").Value = "Fnu Balan"
$ws.Range("D10").Value = "JKR"
$ws.Range("E10").Value = "ENGLISH"
$ws.Range("F10").Value = "202 Hovis Rd Hovis Rd 28164 Nc USA"
$ws.Range("G10").Value = 7046528713

# Row 11
$ws.Range("B11").Value = 45154
$ws.Range("C11").Value = "Amanda Father - Vazquez"
$ws.Range("D11").Value = "GG"
$ws.Range("E11").Value = "ENGLISH"
$ws.Range("F11").Value = "116 cypress Vallejo Ca 94590 United States"
$ws.Range("G11").Value = 70765674356

# Row 12
$ws.Range("B12").Value = 45154
$ws.Range("C12").Value = "Thomas Schenck"
$ws.Range("D12").Value = "GG"
$ws.Range("E12").Value = "ENGLISH"
$ws.Range("F12").Value = "408 210 5th avenue south Saint Petersburg 33701 Florida United States"
$ws.Range("G12").Value = 7273310529

# Row 13
$ws.Range("B13").Value = 45154
$ws.Range("C13").Value = "Mohan Passi"
$ws.Range("D13").Value = "JKR"
$ws.Range("E13").Value = "English"
$ws.Range("F13").Value = "4631 Gresham Drive, Eldorado Hills, 95762, CA, USA"
$ws.Range("G13").Value = 9169330761
$ws.Range("J13").Value = "Mohan Passi 4631 Gresham Drive Eldorado Hills 95762 CA USA 9169330761 Way of living (English) Bangladesh"
$ws.Range("K13").Value = "USA"

# Row 14
$ws.Range("B14").Value = 45154
$ws.Range("C14").Value = "Anantha Bass"
$ws.Range("D14").Value = "JKR"
$ws.Range("E14").Value = "English"
$ws.Range("F14").Value = "202 10404 Salvia Street, Charlotte, 28277, North Carolina, USA"
$ws.Range("G14").Value = 9804285429
$ws.Range("J14").Value = "Anantha Bass 202 10404 salvia street Charlotte 28277 North Carolina United states 9804285429 Way of living (English) UNITED STATES OF AMERICA"
$ws.Range("K14").Value = "USA"

# Row 15
$ws.Range("B15").Value = 45154
$ws.Range("C15").Value = "Thyagarajan Iyer"
$ws.Range("D15").Value = "JKR"
$ws.Range("E15").Value = "English"
$ws.Range("F15").Value = "6092 Elmbridge Dr, San Jose, 95129, CA, USA"
$ws.Range("G15").Value = 9940103805
$ws.Range("J15").Value = "Thyagarajan Iyer 6092 Elmbridge Dr San Jose 95129 Ca USA 9940103805 WhatsApp Way of living (English) Bangladesh"
$ws.Range("K15").Value = "USA"

# Match the author's final cursor position (was D2, now B8)
$ws.Range("B8").Select()
